$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-11-11"
$ws.Range("A12").Value = "November (through 11-11)"

$ws.Range("B12").Value = 13
$ws.Range("C12").Value = 26
$ws.Range("D12").Value = 42
$ws.Range("F12").Value = 19
$ws.Range("G12").Value = 70
$ws.Range("H12").Value = 74

$ws.Range("B13").Value = 271
$ws.Range("C13").Value = 512
$ws.Range("D13").Value = 752
$ws.Range("F13").Value = 501
$ws.Range("G13").Value = 1127
$ws.Range("H13").Value = 1518
